$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.816.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").Value = "'3.765.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'626.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.81%  "

$ws.Range("D6").Value = "'164.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "

$ws.Range("D7").Value = "'3.763.82"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("E11").Value = "  +2.17%  "

$ws.Range("D12").Value = "'6.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").Value = "'35.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "'4.400.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").Value = "'3.732.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "

$ws.Range("D17").Value = "'68.787.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").Value = "'17.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("E19").Value = "  -1.16%  "

$ws.Range("D20").Value = "'7.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").Value = "'464.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("D22").Value = "'9.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("D24").Value = "'82.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("E25").Value = "  -1.38%  "

$ws.Range("D26").Value = "'11.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.96%  "

$ws.Range("D27").Value = "'2.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.58%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "'10.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.83%  "

$ws.Range("D30").Value = "'3.916.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("E32").Value = "  +2.37%  "

$ws.Range("D33").Value = "'7.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.16%  "

$ws.Range("D34").Value = "'28.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.22%  "

$ws.Range("D35").Value = "'0.170"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.07%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "'3.717.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").Value = "'8.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("D41").Value = "'5.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").Value = "'155.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("D46").Value = "'43.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "

$ws.Range("D47").Value = "'46.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").Value = "'0.294"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'1.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.38%  "

$ws.Range("D50").Value = "'8.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("E51").Value = "  -1.21%  "
